$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows before the "LifeRelEmp" row (row 199) for the
# newly added L7206Cust / L7206Emp / L7206Manager tables.
$ws.Rows("199:201").Insert()

$category = "L7-介接外部系統"
$basePath = "[\\192.168.10.16\St1Share(NAS)\SKL\DB\GenTables\L7-介接外部系統\"

$ws.Range("A199").Value = $category
$ws.Range("B199").Value = "L7206Cust"
$ws.Range("C199").Value = "利害關係人借款人檔"
$ws.Range("D199").Formula = '=HYPERLINK("' + $basePath + 'L7206Cust.xlsx]DBD!A1", "連結")'
$ws.Range("E199").Value = "2023年09月22日 14:38:20"

$ws.Range("A200").Value = $category
$ws.Range("B200").Value = "L7206Emp"
$ws.Range("C200").Value = "利害關係人員工檔"
$ws.Range("D200").Formula = '=HYPERLINK("' + $basePath + 'L7206Emp.xlsx]DBD!A1", "連結")'
$ws.Range("E200").Value = "2023年09月22日 14:38:17"

$ws.Range("A201").Value = $category
$ws.Range("B201").Value = "L7206Manager"
$ws.Range("C201").Value = "利害關係人負責人檔"
$ws.Range("D201").Formula = '=HYPERLINK("' + $basePath + 'L7206Manager.xlsx]DBD!A1", "連結")'
$ws.Range("E201").Value = "2023年09月22日 14:38:15"

# Update "last modified" timestamps for several existing tables.
# (Row numbers below are post-insert positions: the three rows inserted
# above shift everything from the old row 199 onward down by three.)
$ws.Range("E8").Value = "2023年09月21日 11:37:51"    # CustMain
$ws.Range("E55").Value = "2023年09月20日 14:04:49"   # Guarantor
$ws.Range("E73").Value = "2023年09月20日 15:31:11"   # BankRemit
$ws.Range("E74").Value = "2023年09月20日 15:55:32"   # BankRmtf
$ws.Range("E113").Value = "2023年09月15日 17:17:54"  # NegAppr01
$ws.Range("E124").Value = "2023年09月20日 11:22:36"  # PfCoOfficerLog
$ws.Range("E137").Value = "2023年09月22日 11:16:56"  # AcDetail
$ws.Range("E163").Value = "2023年09月22日 16:21:41"  # CdComm
$ws.Range("E317").Value = "2023年09月19日 13:03:24"  # MonthlyFacBal (was E314)
$ws.Range("E325").Value = "2023年09月21日 16:27:12"  # MonthlyLM052Loss (was E322)
$ws.Range("E327").Value = "2023年09月18日 14:27:34"  # MonthlyLM055AssetLoss (was E324)
